# Update generated output values in the "展览" (sheet1) and "全部类型" (sheet4) sheets
# to match the refreshed data snapshot (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet1) column F updates ---
$exhibitUpdates = @{
    5  = 1615
    6  = 3224
    7  = 761
    8  = 1971
    9  = 1892
    10 = 968
    13 = 1569
    14 = 333
    17 = 1383
    18 = 477
    19 = 589
    20 = 279
    21 = 10387
    22 = 9573
    23 = 829
    24 = 632
    25 = 1799
    26 = 128
    27 = 360
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# --- 全部类型 (sheet4) column F updates ---
$allUpdates = @{
    7  = 1615
    8  = 3224
    9  = 761
    10 = 1971
    11 = 1892
    12 = 968
    15 = 1569
    16 = 333
    21 = 1383
    22 = 477
    23 = 589
    24 = 279
    25 = 10387
    26 = 9574
    27 = 829
    28 = 632
    29 = 1799
    32 = 128
    33 = 360
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
